$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to retain text storage (values like "1.002" would
# otherwise be auto-converted to numbers by the Value setter); ClearFormats
# afterwards restores the original (unstyled) appearance.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.052.95"
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").Value = "1.890.35"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "313.88"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.5023"
$ws.Range("E7").Value = "  -4.69%  "
$ws.Range("D8").Value = "0.3894"
$ws.Range("E8").Value = "  -1.76%  "
$ws.Range("D9").Value = "0.09204"
$ws.Range("E9").Value = "  -4.85%  "
$ws.Range("D10").Value = "1.129"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").Value = "41.78"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "6.371"
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("D13").Value = "20.81"
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("D14").Value = "1.890.94"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("D15").Value = "7.302"
$ws.Range("E15").Value = "  -3.39%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "0.00001108"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("D18").Value = "91.80"
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("D19").Value = "0.06628"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "17.85"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "6.232"
$ws.Range("E22").Value = "  -1.61%  "
$ws.Range("D23").Value = "28.138.20"
$ws.Range("E23").Value = "  -1.74%  "
$ws.Range("D24").Value = "11.37"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").Value = "2.321"
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("D26").Value = "2.114.16"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").Value = "2.554"
$ws.Range("E27").Value = "  -5.19%  "
$ws.Range("D28").Value = "158.44"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").Value = "20.78"
$ws.Range("E29").Value = "  -2.22%  "
$ws.Range("D30").Value = "126.92"
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("D31").Value = "1.074"
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("E32").Value = "  -3.12%  "
$ws.Range("D33").Value = "5.593"
$ws.Range("E33").Value = "  -2.75%  "
$ws.Range("D34").Value = "3.603"
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("D35").Value = "9.494"
$ws.Range("E35").Value = "  -3.86%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "1.348"
$ws.Range("E36").Value = "  +13.41%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.06583"
$ws.Range("E37").Value = "  -2.25%  "
$ws.Range("D38").Value = "0.02415"
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("D39").Value = "0.2197"
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("D40").Value = "1.216"
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("D41").Value = "0.6455"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").Value = "11.57"
$ws.Range("E42").Value = "  -1.94%  "
$ws.Range("D43").Value = "4.951"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "13.43"
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("D46").Value = "0.6059"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D48").Value = "3.691"
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("D49").Value = "2.002"
$ws.Range("E49").Value = "  -1.56%  "
$ws.Range("D50").Value = "121.75"
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("D51").Value = "1.195"
$ws.Range("E51").Value = "  -1.22%  "

$priceRange.ClearFormats()

"done"
